$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Owners")

$ws.Range("D37").Value = "NO FB PAGE"
$ws.Range("D46").Value = "info@pleasantonsmiles.com"
$ws.Range("E46").Value = "(925) 462-1990"
$ws.Range("D49").Value = "No email found"
$ws.Range("E49").Value = "925-433-7809"
$ws.Range("D53").Value = "No email found"
$ws.Range("E53").Value = "'9256009888"
$ws.Range("D55").Value = "No email found"
$ws.Range("E55").Value = "(925) 249-1130"
$ws.Range("D57").Value = "No email found"
$ws.Range("E57").Value = "'214872808"
$ws.Range("D58").Value = "NO FB PAGE"
$ws.Range("E58").Value = "'9256009888"
$ws.Range("D59").Value = "NO FB PAGE"
$ws.Range("E59").Value = "(925) 485-1172"
$ws.Range("D61").Value = "NO FB PAGE"
$ws.Range("E61").Value = "NO FB PAGE"
$ws.Range("D62").Value = "contact@pearlsandteeth.com"
$ws.Range("E62").Value = "925-523-3864"
$ws.Range("D63").Value = "No email found"
$ws.Range("E63").Value = "(925) 846-6308"
$ws.Range("D64").Value = "NO FB PAGE"
$ws.Range("E64").Value = "925-425-7545"
$ws.Range("D65").Value = "arnoldjacobsstaff@att.net"
$ws.Range("E65").Value = "(925)846-3968"
$ws.Range("D66").Value = "info@eastbaydental.com"
$ws.Range("E66").Value = "(510) 818-9000"
$ws.Range("D76").Value = "normanrwong@ymail.com"
$ws.Range("E76").Value = "(925) 846-5506"
$ws.Range("D77").Value = "info@redwoodcitydental.com"
$ws.Range("E77").Value = "(925) 551-6464"
$ws.Range("D79").Value = "kjudson321@gmail.com"
$ws.Range("E79").Value = "925-462-1990"
$ws.Range("D80").Value = "dublinranchdental@yahoo.com"
$ws.Range("E80").Value = "(925) 999-9088"
$ws.Range("D81").Value = "pleasantondds@gmail.com"
$ws.Range("E81").Value = "'9256009888"
$ws.Range("D82").Value = "No email found"
$ws.Range("E82").Value = "(415) 380-3600"
$ws.Range("D83").Value = "syoondds@gmail.com"
$ws.Range("E83").Value = "(707) 422-7003"
$ws.Range("D85").Value = "fairfielddentists@smilegeneration.com"
$ws.Range("E85").Value = "707-399-9082"
$ws.Range("D86").Value = "No email found"
$ws.Range("E86").Value = "(707) 716-1715"
$ws.Range("D92").Value = "galvandentalcorporation@gmail.com"
$ws.Range("E92").Value = "(925)676-6363"
$ws.Range("D93").Value = "victodonto@msn.com"
$ws.Range("E93").Value = "99931-7003"
$ws.Range("D94").Value = "No email found"
$ws.Range("E94").Value = "(409) 883-5300"
$ws.Range("D96").Value = "NO FB PAGE"
$ws.Range("E96").Value = "NO FB PAGE"
$ws.Range("D111").Value = "No email found"
$ws.Range("E111").Value = "916-246-1502"
$ws.Range("D114").Value = "No email found"
$ws.Range("E114").Value = "(916) 723-3368"
$ws.Range("D121").Value = "dentistsofoldtorrance@smilegeneration.com"
$ws.Range("E121").Value = "(424) 320-6584"
$ws.Range("D123").Value = "mehr@sparklefamilydentist.com"
$ws.Range("E123").Value = "424-378-1168"
$ws.Range("D130").Value = "clientservice@yourhealthcontact.com"
$ws.Range("E130").Value = "'3633090"
$ws.Range("D148").Value = "No email found"
$ws.Range("E148").Value = "(310) 782-2008"
$ws.Range("D149").Value = "No email found"
$ws.Range("E149").Value = "(310) 328-9700"
$ws.Range("D150").Value = "No email found"
$ws.Range("E150").Value = "310-534-3002"
$ws.Range("D154").Value = "No email found"
$ws.Range("E154").Value = "(310) 483-7779"
$ws.Range("D156").Value = "support@dentalinsider.com"
$ws.Range("E156").Value = "(310) 320-3264"
$ws.Range("D157").Value = "No email found"
$ws.Range("E157").Value = "(310) 543-1655"
$ws.Range("D162").ClearContents()
$ws.Range("E162").ClearContents()
$ws.Range("D163").Value = "louis_yang2@yahoo.com"
$ws.Range("E163").Value = "94587-2743"
$ws.Range("D168").Value = "smile@igasakidental.com"
$ws.Range("E168").Value = "NO FB PAGE"
$ws.Range("D169").Value = "IPIofTorrance@gmail.com"
$ws.Range("E169").Value = "(310)320-5661"
$ws.Range("D173").Value = "NO FB PAGE"
$ws.Range("E173").Value = "NO FB PAGE"
$ws.Range("D174").Value = "WCDCustomerService@westcoastdental.com"
$ws.Range("E174").Value = "(888) 329-8111"
$ws.Range("D176").Value = "No email found"
$ws.Range("E176").Value = "(669) 209-9319"
$ws.Range("D177").ClearContents()
$ws.Range("E177").ClearContents()
$ws.Range("D178").ClearContents()
$ws.Range("E178").ClearContents()
$ws.Range("D184").Value = "info@nohosmilecenter.com"
$ws.Range("E184").Value = "(818) 505-0106"
$ws.Range("D192").Value = "info@bookdok.com"
$ws.Range("E192").Value = "818-788-2023"
$ws.Range("D195").Value = "No email found"
$ws.Range("E195").Value = "(559) 325-0700"
$ws.Range("D196").Value = "dennisdentalcare@gmail.com"
$ws.Range("E196").Value = "(559) 570-6981"
$ws.Range("D204").Value = "plazadentalla@gmail.com"
$ws.Range("E204").Value = "(818) 426-6654"
$ws.Range("D207").Value = "No email found"
$ws.Range("E207").Value = "213-680-2808"
$ws.Range("D211").Value = "No email found"
$ws.Range("E211").Value = "(213) 534-6856"
$ws.Range("D216").Value = "ajg@u.arizona.edu"
$ws.Range("E216").Value = "91302-3014"
$ws.Range("D220").Value = "NO FB PAGE"
$ws.Range("E220").Value = "'0506627"
$ws.Range("D222").Value = "No email found"
$ws.Range("E222").Value = "No phone found"
$ws.Range("D227").Value = "drmahdavi@skylinedental.com"
$ws.Range("E227").Value = "(805) 212-5091"

# Update C168, C185, C192 literal padding text (trailing underscore-escaped CRs)
$ws.Range("C168").Value = "Dr. Alan Igasaki_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
Dr. Howard Igasaki"
$ws.Range("C185").Value = "Mary Yazdan, DDS_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
Dr. Hameed Nejat DMD"
$ws.Range("C192").Value = "Robert Tingillian, DDS_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
Tom Shanakian, DDS"
